# Generate Report for Handback
# ------------------------------------------------------------------
# This mirrors the "handback" report generation step: the previous
# "Ready for handoff" status is flipped to "Handed back: in sync with
# en-US", the handback timestamp columns (H) get a real datetime, and
# two new columns are populated per language sheet:
#   F = Latest Target File   (same file identity as the source .md, col A)
#   G = Latest Handback File (same file identity as the handoff .xlf, col D)
# with matching hyperlinks / hyperlink styling.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# ------------------------------------------------------------------
# 1. Flip every "Ready for handoff" status cell to the handback status.
# ------------------------------------------------------------------
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus

$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# ------------------------------------------------------------------
# 2. Latest Handback DateTime (column H) now has a real timestamp.
# ------------------------------------------------------------------
$wsZhCn.Range("H2").Value = "2016-03-11 20:43:45"
$wsZhCn.Range("H3").Value = "2016-03-11 20:43:45"

$wsDeDe.Range("H2").Value = "2016-03-11 20:43:50"
$wsDeDe.Range("H3").Value = "2016-03-11 20:43:50"

# ------------------------------------------------------------------
# 3. Populate "Latest Target File" (F) / "Latest Handback File" (G)
#    for the zh-cn sheet, with hyperlinks that mirror columns A / D.
# ------------------------------------------------------------------
$linkColor = 15570276   # OLE BGR for RGB(100,149,237) / hex 6495ED - matches the workbook's HyperLink style

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/daa119de4003f3d1732debde1fea371bd930a7d9/e2e/21df1751-2a26-4c19-8679-12b22d725b86.md", "", "", "21df1751-2a26-4c19-8679-12b22d725b86.md")
$wsZhCn.Range("F2").Font.Underline = 2
$wsZhCn.Range("F2").Font.Color = $linkColor

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5dddf3e89bfdd2b7cdc0b38c89b17fb0e0e1d036/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/21df1751-2a26-4c19-8679-12b22d725b86.e32f2d9a09563552f6b45587a56d141ec9681cec.zh-cn.xlf", "", "", "21df1751-2a26-4c19-8679-12b22d725b86.e32f2d9a09563552f6b45587a56d141ec9681cec.zh-cn.xlf")
$wsZhCn.Range("G2").Font.Underline = 2
$wsZhCn.Range("G2").Font.Color = $linkColor

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/daa119de4003f3d1732debde1fea371bd930a7d9/e2e/582cfb1a-645f-41e2-a5ff-9db963d3d27a.md", "", "", "582cfb1a-645f-41e2-a5ff-9db963d3d27a.md")
$wsZhCn.Range("F3").Font.Underline = 2
$wsZhCn.Range("F3").Font.Color = $linkColor

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5dddf3e89bfdd2b7cdc0b38c89b17fb0e0e1d036/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/582cfb1a-645f-41e2-a5ff-9db963d3d27a.8f19445d09d6c0fd54db8a6edc35223fdb534180.zh-cn.xlf", "", "", "582cfb1a-645f-41e2-a5ff-9db963d3d27a.8f19445d09d6c0fd54db8a6edc35223fdb534180.zh-cn.xlf")
$wsZhCn.Range("G3").Font.Underline = 2
$wsZhCn.Range("G3").Font.Color = $linkColor

# ------------------------------------------------------------------
# 4. Same for the de-de sheet.
# ------------------------------------------------------------------
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/daa119de4003f3d1732debde1fea371bd930a7d9/e2e/21df1751-2a26-4c19-8679-12b22d725b86.md", "", "", "21df1751-2a26-4c19-8679-12b22d725b86.md")
$wsDeDe.Range("F2").Font.Underline = 2
$wsDeDe.Range("F2").Font.Color = $linkColor

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5d54a7a3ff747f5d08b9d2b7577274501bc17809/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/21df1751-2a26-4c19-8679-12b22d725b86.e32f2d9a09563552f6b45587a56d141ec9681cec.de-de.xlf", "", "", "21df1751-2a26-4c19-8679-12b22d725b86.e32f2d9a09563552f6b45587a56d141ec9681cec.de-de.xlf")
$wsDeDe.Range("G2").Font.Underline = 2
$wsDeDe.Range("G2").Font.Color = $linkColor

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/daa119de4003f3d1732debde1fea371bd930a7d9/e2e/582cfb1a-645f-41e2-a5ff-9db963d3d27a.md", "", "", "582cfb1a-645f-41e2-a5ff-9db963d3d27a.md")
$wsDeDe.Range("F3").Font.Underline = 2
$wsDeDe.Range("F3").Font.Color = $linkColor

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5d54a7a3ff747f5d08b9d2b7577274501bc17809/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/582cfb1a-645f-41e2-a5ff-9db963d3d27a.8f19445d09d6c0fd54db8a6edc35223fdb534180.de-de.xlf", "", "", "582cfb1a-645f-41e2-a5ff-9db963d3d27a.8f19445d09d6c0fd54db8a6edc35223fdb534180.de-de.xlf")
$wsDeDe.Range("G3").Font.Underline = 2
$wsDeDe.Range("G3").Font.Color = $linkColor

Write-Output "Handback report generated."
